$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.841.76'
$ws.Range("E2").Value = '  +3.16%  '

$ws.Range("D3").Value = '1.879.25'
$ws.Range("E3").Value = '  +3.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.63'
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("E7").Value = '  +1.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3952'
$ws.Range("E8").Value = '  +3.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07926'
$ws.Range("E9").Value = '  +1.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9776'
$ws.Range("E10").Value = '  +2.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.44'
$ws.Range("E11").Value = '  +2.81%  '

$ws.Range("D12").Value = '1.915.74'
$ws.Range("E12").Value = '  +8.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.763'
$ws.Range("E13").Value = '  +2.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.034'
$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06989'
$ws.Range("E15").Value = '  +1.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.74'
$ws.Range("E16").Value = '  +2.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  +0.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001011'
$ws.Range("E18").Value = '  +1.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.04'
$ws.Range("E19").Value = '  +1.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").Value = '28.851.03'
$ws.Range("E21").Value = '  +3.09%  '

$ws.Range("E22").Value = '  +1.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.13'
$ws.Range("E23").Value = '  +2.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.121'
$ws.Range("E24").Value = '  -0.17%  '

$ws.Range("D25").Value = '2.072.62'
$ws.Range("E25").Value = '  +4.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.75'
$ws.Range("E26").Value = '  +1.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.41'
$ws.Range("E27").Value = '  +1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.775'
$ws.Range("E28").Value = '  +1.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.010'
$ws.Range("E29").Value = '  +2.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.86'
$ws.Range("E30").Value = '  +3.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09408'
$ws.Range("E31").Value = '  +1.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9433'
$ws.Range("E32").Value = '  +1.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.331'
$ws.Range("E33").Value = '  +0.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.355'
$ws.Range("E34").Value = '  +3.91%  '

$ws.Range("E35").Value = '  -2.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05932'
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02133'
$ws.Range("E37").Value = '  -0.40%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.151'
$ws.Range("E38").Value = '  +0.62%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.945'
$ws.Range("E39").Value = '  +5.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5741'
$ws.Range("E40").Value = '  +3.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1796'
$ws.Range("E41").Value = '  +1.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.988'
$ws.Range("E42").Value = '  +1.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07249'
$ws.Range("E43").Value = '  +3.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.83'
$ws.Range("E44").Value = '  +1.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5353'
$ws.Range("E45").Value = '  +2.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.160'
$ws.Range("E46").Value = '  -4.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.141'
$ws.Range("E47").Value = '  -4.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.855'
$ws.Range("E48").Value = '  +1.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.43'
$ws.Range("E49").Value = '  +1.92%  '

$ws.Range("E50").Value = '  +3.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.006'
$ws.Range("E51").Value = '  +0.32%  '
